$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CompStat")

# --- Header text (Volume number + date range) ---
$ws.Range("A8").Value = "Volume 33   Number  2"
$ws.Range("C9").Value = "Report Covering the Week  1/5/2026  Through  1/11/2026"

# --- Cells changing type/style (number <-> dash/undefined placeholder text) ---
# Donor cells (untouched elsewhere) supply exact style + content via copy/paste-special
#   C14 -> style 13, text "0"      (s=13 t=s v=20)
#   E14 -> style 13, text "***.*"  (s=13 t=s v=21)
#   D15 -> style 15, numeric
#   E15 -> style 14, numeric

$ws.Range("C14").Copy()
$ws.Range("C16").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C16").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)

$ws.Range("E15").Copy()
$ws.Range("M16").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("M16").Value = -75

$ws.Range("D15").Copy()
$ws.Range("C18").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C18").Value = 1

$ws.Range("C14").Copy()
$ws.Range("D18").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("D18").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)

$ws.Range("E14").Copy()
$ws.Range("E18").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E18").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)

$ws.Range("D15").Copy()
$ws.Range("I18").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("I18").Value = 1

$ws.Range("E15").Copy()
$ws.Range("L18").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("L18").Value = 0

$ws.Range("E15").Copy()
$ws.Range("M18").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("M18").Value = -50

$ws.Range("E15").Copy()
$ws.Range("M19").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("M19").Value = -14.285714285714

$ws.Range("D15").Copy()
$ws.Range("C20").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("C20").Value = 1

$ws.Range("D15").Copy()
$ws.Range("I20").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("I20").Value = 1

$ws.Range("E15").Copy()
$ws.Range("L22").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("L22").Value = -100

$ws.Range("C14").Copy()
$ws.Range("D25").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("D25").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)

$ws.Range("E14").Copy()
$ws.Range("E25").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E25").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)

$ws.Range("D15").Copy()
$ws.Range("I25").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("I25").Value = 2

$ws.Range("C14").Copy()
$ws.Range("G31").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("G31").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)

$ws.Range("E14").Copy()
$ws.Range("H31").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("H31").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues)

# --- Simple value-only updates (style/type unchanged) ---
$ws.Range("G15").Value = 2
$ws.Range("J15").Value = 2
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -100
$ws.Range("F16").Value = 3
$ws.Range("H16").Value = -57.142857142857
$ws.Range("J16").Value = 5
$ws.Range("K16").Value = -80
$ws.Range("L16").Value = -83.333333333333
$ws.Range("N16").Value = -92.857142857142
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 7
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = -36.363636363636
$ws.Range("I17").Value = 4
$ws.Range("J17").Value = 5
$ws.Range("K17").Value = -20
$ws.Range("L17").Value = 33.333333333333
$ws.Range("M17").Value = 300
$ws.Range("F18").Value = 3
$ws.Range("H18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("N18").Value = -96.296296296296
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 26
$ws.Range("G19").Value = 34
$ws.Range("H19").Value = -23.529411764705
$ws.Range("I19").Value = 6
$ws.Range("J19").Value = 9
$ws.Range("K19").Value = -33.333333333333
$ws.Range("L19").Value = 0
$ws.Range("N19").Value = -66.666666666666
$ws.Range("F20").Value = 4
$ws.Range("N20").Value = -90.909090909090
$ws.Range("C21").Value = 8
$ws.Range("D21").Value = 9
$ws.Range("E21").Value = -11.111111111111
$ws.Range("F21").Value = 43
$ws.Range("G21").Value = 57
$ws.Range("H21").Value = -24.561403508771
$ws.Range("I21").Value = 13
$ws.Range("J21").Value = 22
$ws.Range("K21").Value = -40.909090909090
$ws.Range("L21").Value = -18.75
$ws.Range("M21").Value = -18.75
$ws.Range("N21").Value = -83.333333333333
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -66.666666666666
$ws.Range("J22").Value = 2
$ws.Range("C23").Value = 1
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 0
$ws.Range("G23").Value = 9
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 3
$ws.Range("J23").Value = 3
$ws.Range("L23").Value = -25
$ws.Range("M23").Value = -25
$ws.Range("C24").Value = 9
$ws.Range("D24").Value = 5
$ws.Range("E24").Value = 80
$ws.Range("F24").Value = 37
$ws.Range("G24").Value = 22
$ws.Range("H24").Value = 68.181818181818
$ws.Range("I24").Value = 10
$ws.Range("J24").Value = 9
$ws.Range("K24").Value = 11.111111111111
$ws.Range("L24").Value = 100
$ws.Range("M24").Value = -37.5
$ws.Range("C25").Value = 2
$ws.Range("F25").Value = 11
$ws.Range("G25").Value = 3
$ws.Range("H25").Value = 266.666666666667
$ws.Range("K25").Value = 100
$ws.Range("L25").Value = 0
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 6
$ws.Range("E26").Value = -50
$ws.Range("F26").Value = 8
$ws.Range("G26").Value = 26
$ws.Range("H26").Value = -69.230769230769
$ws.Range("I26").Value = 4
$ws.Range("J26").Value = 12
$ws.Range("K26").Value = -66.666666666666
$ws.Range("M26").Value = -20
$ws.Range("D27").Value = 2
$ws.Range("G27").Value = 3
$ws.Range("J27").Value = 3
$ws.Range("C28").Value = 4
$ws.Range("F28").Value = 8
$ws.Range("I28").Value = 5
$ws.Range("J45").Value = 48
$ws.Range("K45").Value = -63.909774436090
$ws.Range("L45").Value = -57.894736842105
$ws.Range("M45").Value = -86.666666666666
$ws.Range("N45").Value = -92.452830188679
$ws.Range("J46").Value = 760
$ws.Range("K46").Value = -19.576719576719
$ws.Range("L46").Value = -42.727957799547
$ws.Range("M46").Value = -72.353583121135
$ws.Range("N46").Value = -77.514792899408

$excel.CutCopyMode = $false
